# edit.ps1 - RPA datasets push 2023-11-04
# Rewrites the full data table (rows 2-27) on Sheet1 to match the
# updated IPO underwriting dataset (new deals added, stale deals removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force plain text so date-like strings ("2023-09-21") are not
    # auto-converted into Excel date serials, then drop the temporary
    # text number-format so the cell keeps the workbook default style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2: 두산로보틱스
Set-TextCell $ws.Range("A2") "CS"
Set-TextCell $ws.Range("B2") "2023-09-21"
Set-TextCell $ws.Range("C2") "두산로보틱스"
Set-TextCell $ws.Range("D2") "한국, 미래"
Set-TextCell $ws.Range("E2") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F2") "2023-09-26"
Set-TextCell $ws.Range("G2") "2023-10-05"
$ws.Range("H2").Value = 42120
$ws.Range("I2").Value = 16200000
$ws.Range("J2").Value = 26000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 10

# Row 3: KB제27호스팩
Set-TextCell $ws.Range("A3") "KB"
Set-TextCell $ws.Range("B3") "2023-10-24"
Set-TextCell $ws.Range("C3") "KB제27호스팩"
Set-TextCell $ws.Range("D3") "KB"
Set-TextCell $ws.Range("E3") "KB"
Set-TextCell $ws.Range("F3") "2023-10-27"
Set-TextCell $ws.Range("G3") "2023-11-03"
$ws.Range("H3").Value = 25000
$ws.Range("I3").Value = 12500000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100

# Row 4: 한싹
Set-TextCell $ws.Range("A4") "KB"
Set-TextCell $ws.Range("B4") "2023-09-19"
Set-TextCell $ws.Range("C4") "한싹"
Set-TextCell $ws.Range("D4") "KB"
Set-TextCell $ws.Range("E4") "KB"
Set-TextCell $ws.Range("F4") "2023-09-22"
Set-TextCell $ws.Range("G4") "2023-10-04"
$ws.Range("H4").Value = 18750
$ws.Range("I4").Value = 1500000
$ws.Range("J4").Value = 12500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100

# Row 5: 두산로보틱스
Set-TextCell $ws.Range("A5") "KB"
Set-TextCell $ws.Range("B5") "2023-09-21"
Set-TextCell $ws.Range("C5") "두산로보틱스"
Set-TextCell $ws.Range("D5") "한국, 미래"
Set-TextCell $ws.Range("E5") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F5") "2023-09-26"
Set-TextCell $ws.Range("G5") "2023-10-05"
$ws.Range("H5").Value = 42120
$ws.Range("I5").Value = 16200000
$ws.Range("J5").Value = 26000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 10

# Row 6: 두산로보틱스
Set-TextCell $ws.Range("A6") "NH"
Set-TextCell $ws.Range("B6") "2023-09-21"
Set-TextCell $ws.Range("C6") "두산로보틱스"
Set-TextCell $ws.Range("D6") "한국, 미래"
Set-TextCell $ws.Range("E6") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F6") "2023-09-26"
Set-TextCell $ws.Range("G6") "2023-10-05"
$ws.Range("H6").Value = 42120
$ws.Range("I6").Value = 16200000
$ws.Range("J6").Value = 26000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 10

# Row 7: 유진테크놀로지
Set-TextCell $ws.Range("A7") "NH"
Set-TextCell $ws.Range("B7") "2023-10-23"
Set-TextCell $ws.Range("C7") "유진테크놀로지"
Set-TextCell $ws.Range("D7") "NH"
Set-TextCell $ws.Range("E7") "NH"
Set-TextCell $ws.Range("F7") "2023-10-26"
Set-TextCell $ws.Range("G7") "2023-11-02"
$ws.Range("H7").Value = 17841.194
$ws.Range("I7").Value = 1049482
$ws.Range("J7").Value = 17000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 100

# Row 8: 밀리의서재
Set-TextCell $ws.Range("A8") "미래"
Set-TextCell $ws.Range("B8") "2023-09-18"
Set-TextCell $ws.Range("C8") "밀리의서재"
Set-TextCell $ws.Range("D8") "미래"
Set-TextCell $ws.Range("E8") "미래"
Set-TextCell $ws.Range("F8") "2023-09-21"
Set-TextCell $ws.Range("G8") "2023-09-27"
$ws.Range("H8").Value = 34500
$ws.Range("I8").Value = 1500000
$ws.Range("J8").Value = 23000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 100

# Row 9: 신성에스티
Set-TextCell $ws.Range("A9") "미래"
Set-TextCell $ws.Range("B9") "2023-10-10"
Set-TextCell $ws.Range("C9") "신성에스티"
Set-TextCell $ws.Range("D9") "미래"
Set-TextCell $ws.Range("E9") "미래"
Set-TextCell $ws.Range("F9") "2023-10-13"
Set-TextCell $ws.Range("G9") "2023-10-19"
$ws.Range("H9").Value = 52000
$ws.Range("I9").Value = 2000000
$ws.Range("J9").Value = 26000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100

# Row 10: 퓨릿
Set-TextCell $ws.Range("A10") "미래"
Set-TextCell $ws.Range("B10") "2023-10-05"
Set-TextCell $ws.Range("C10") "퓨릿"
Set-TextCell $ws.Range("D10") "미래"
Set-TextCell $ws.Range("E10") "미래"
Set-TextCell $ws.Range("F10") "2023-10-11"
Set-TextCell $ws.Range("G10") "2023-10-18"
$ws.Range("H10").Value = 44265.9
$ws.Range("I10").Value = 4137000
$ws.Range("J10").Value = 10700
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 100

# Row 11: 두산로보틱스
Set-TextCell $ws.Range("A11") "미래"
Set-TextCell $ws.Range("B11") "2023-09-21"
Set-TextCell $ws.Range("C11") "두산로보틱스"
Set-TextCell $ws.Range("D11") "한국, 미래"
Set-TextCell $ws.Range("E11") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F11") "2023-09-26"
Set-TextCell $ws.Range("G11") "2023-10-05"
$ws.Range("H11").Value = 126360
$ws.Range("I11").Value = 16200000
$ws.Range("J11").Value = 26000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 30

# Row 12: 레뷰코퍼레이션
Set-TextCell $ws.Range("A12") "삼성"
Set-TextCell $ws.Range("B12") "2023-09-19"
Set-TextCell $ws.Range("C12") "레뷰코퍼레이션"
Set-TextCell $ws.Range("D12") "삼성"
Set-TextCell $ws.Range("E12") "삼성"
Set-TextCell $ws.Range("F12") "2023-09-22"
Set-TextCell $ws.Range("G12") "2023-10-06"
$ws.Range("H12").Value = 33600
$ws.Range("I12").Value = 2240000
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100

# Row 13: 상상인제4호스팩
Set-TextCell $ws.Range("A13") "상상인"
Set-TextCell $ws.Range("B13") "2023-09-04"
Set-TextCell $ws.Range("C13") "상상인제4호스팩"
Set-TextCell $ws.Range("D13") "상상인"
Set-TextCell $ws.Range("E13") "상상인"
Set-TextCell $ws.Range("F13") "2023-09-07"
Set-TextCell $ws.Range("G13") "2023-09-14"
$ws.Range("H13").Value = 9000
$ws.Range("I13").Value = 4500000
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100

# Row 14: 인스웨이브시스템즈
Set-TextCell $ws.Range("A14") "신영"
Set-TextCell $ws.Range("B14") "2023-09-14"
Set-TextCell $ws.Range("C14") "인스웨이브시스템즈"
Set-TextCell $ws.Range("D14") "신영"
Set-TextCell $ws.Range("E14") "신영"
Set-TextCell $ws.Range("F14") "2023-09-19"
Set-TextCell $ws.Range("G14") "2023-09-25"
$ws.Range("H14").Value = 26400
$ws.Range("I14").Value = 1100000
$ws.Range("J14").Value = 24000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 100

# Row 15: 두산로보틱스
Set-TextCell $ws.Range("A15") "신영"
Set-TextCell $ws.Range("B15") "2023-09-21"
Set-TextCell $ws.Range("C15") "두산로보틱스"
Set-TextCell $ws.Range("D15") "한국, 미래"
Set-TextCell $ws.Range("E15") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F15") "2023-09-26"
Set-TextCell $ws.Range("G15") "2023-10-05"
$ws.Range("H15").Value = 12636
$ws.Range("I15").Value = 16200000
$ws.Range("J15").Value = 26000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 3

# Row 16: 유투바이오
Set-TextCell $ws.Range("A16") "신한"
Set-TextCell $ws.Range("B16") "2023-10-23"
Set-TextCell $ws.Range("C16") "유투바이오"
Set-TextCell $ws.Range("D16") "신한"
Set-TextCell $ws.Range("E16") "신한"
Set-TextCell $ws.Range("F16") "2023-10-26"
Set-TextCell $ws.Range("G16") "2023-11-02"
$ws.Range("H16").Value = 4966.368
$ws.Range("I16").Value = 1128720
$ws.Range("J16").Value = 4400
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100

# Row 17: 신한제11호스팩
Set-TextCell $ws.Range("A17") "신한"
Set-TextCell $ws.Range("B17") "2023-09-19"
Set-TextCell $ws.Range("C17") "신한제11호스팩"
Set-TextCell $ws.Range("D17") "신한"
Set-TextCell $ws.Range("E17") "신한"
Set-TextCell $ws.Range("F17") "2023-09-22"
Set-TextCell $ws.Range("G17") "2023-10-04"
$ws.Range("H17").Value = 36000
$ws.Range("I17").Value = 18000000
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 100

# Row 18: 두산로보틱스
Set-TextCell $ws.Range("A18") "유비에스"
Set-TextCell $ws.Range("B18") "2023-09-21"
Set-TextCell $ws.Range("C18") "두산로보틱스"
Set-TextCell $ws.Range("D18") "한국, 미래"
Set-TextCell $ws.Range("E18") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F18") "2023-09-26"
Set-TextCell $ws.Range("G18") "2023-10-05"
$ws.Range("H18").Value = 4212
$ws.Range("I18").Value = 16200000
$ws.Range("J18").Value = 26000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 1

# Row 19: 아이엠티
Set-TextCell $ws.Range("A19") "유안타"
Set-TextCell $ws.Range("B19") "2023-09-18"
Set-TextCell $ws.Range("C19") "아이엠티"
Set-TextCell $ws.Range("D19") "유안타"
Set-TextCell $ws.Range("E19") "유안타, 유진"
Set-TextCell $ws.Range("F19") "2023-09-21"
Set-TextCell $ws.Range("G19") "2023-10-10"
$ws.Range("H19").Value = 15484
$ws.Range("I19").Value = 1580000
$ws.Range("J19").Value = 14000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 70

# Row 20: 아이엠티
Set-TextCell $ws.Range("A20") "유진"
Set-TextCell $ws.Range("B20") "2023-09-18"
Set-TextCell $ws.Range("C20") "아이엠티"
Set-TextCell $ws.Range("D20") "유안타"
Set-TextCell $ws.Range("E20") "유안타, 유진"
Set-TextCell $ws.Range("F20") "2023-09-21"
Set-TextCell $ws.Range("G20") "2023-10-10"
$ws.Range("H20").Value = 6636
$ws.Range("I20").Value = 1580000
$ws.Range("J20").Value = 14000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 30

# Row 21: 워트
Set-TextCell $ws.Range("A21") "키움"
Set-TextCell $ws.Range("B21") "2023-10-16"
Set-TextCell $ws.Range("C21") "워트"
Set-TextCell $ws.Range("D21") "키움"
Set-TextCell $ws.Range("E21") "키움"
Set-TextCell $ws.Range("F21") "2023-10-19"
Set-TextCell $ws.Range("G21") "2023-10-26"
$ws.Range("H21").Value = 26000
$ws.Range("I21").Value = 4000000
$ws.Range("J21").Value = 6500
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 100

# Row 22: 두산로보틱스
Set-TextCell $ws.Range("A22") "키움"
Set-TextCell $ws.Range("B22") "2023-09-21"
Set-TextCell $ws.Range("C22") "두산로보틱스"
Set-TextCell $ws.Range("D22") "한국, 미래"
Set-TextCell $ws.Range("E22") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F22") "2023-09-26"
Set-TextCell $ws.Range("G22") "2023-10-05"
$ws.Range("H22").Value = 12636
$ws.Range("I22").Value = 16200000
$ws.Range("J22").Value = 26000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3

# Row 23: 에스엘에스바이오
Set-TextCell $ws.Range("A23") "하나"
Set-TextCell $ws.Range("B23") "2023-10-10"
Set-TextCell $ws.Range("C23") "에스엘에스바이오"
Set-TextCell $ws.Range("D23") "하나"
Set-TextCell $ws.Range("E23") "하나"
Set-TextCell $ws.Range("F23") "2023-10-13"
Set-TextCell $ws.Range("G23") "2023-10-20"
$ws.Range("H23").Value = 5390
$ws.Range("I23").Value = 770000
$ws.Range("J23").Value = 7000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100

# Row 24: 두산로보틱스
Set-TextCell $ws.Range("A24") "하나"
Set-TextCell $ws.Range("B24") "2023-09-21"
Set-TextCell $ws.Range("C24") "두산로보틱스"
Set-TextCell $ws.Range("D24") "한국, 미래"
Set-TextCell $ws.Range("E24") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F24") "2023-09-26"
Set-TextCell $ws.Range("G24") "2023-10-05"
$ws.Range("H24").Value = 12636
$ws.Range("I24").Value = 16200000
$ws.Range("J24").Value = 26000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 3

# Row 25: 퀄리타스반도체
Set-TextCell $ws.Range("A25") "한국"
Set-TextCell $ws.Range("B25") "2023-10-18"
Set-TextCell $ws.Range("C25") "퀄리타스반도체"
Set-TextCell $ws.Range("D25") "한국"
Set-TextCell $ws.Range("E25") "한국"
Set-TextCell $ws.Range("F25") "2023-10-23"
Set-TextCell $ws.Range("G25") "2023-10-27"
$ws.Range("H25").Value = 30600
$ws.Range("I25").Value = 1800000
$ws.Range("J25").Value = 17000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 100

# Row 26: 두산로보틱스
Set-TextCell $ws.Range("A26") "한국"
Set-TextCell $ws.Range("B26") "2023-09-21"
Set-TextCell $ws.Range("C26") "두산로보틱스"
Set-TextCell $ws.Range("D26") "한국, 미래"
Set-TextCell $ws.Range("E26") "한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)"
Set-TextCell $ws.Range("F26") "2023-09-26"
Set-TextCell $ws.Range("G26") "2023-10-05"
$ws.Range("H26").Value = 126360
$ws.Range("I26").Value = 16200000
$ws.Range("J26").Value = 26000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 30

# Row 27: 에이치엠씨제6호스팩
Set-TextCell $ws.Range("A27") "현대차"
Set-TextCell $ws.Range("B27") "2023-09-25"
Set-TextCell $ws.Range("C27") "에이치엠씨제6호스팩"
Set-TextCell $ws.Range("D27") "현대차"
Set-TextCell $ws.Range("E27") "현대차"
Set-TextCell $ws.Range("F27") "2023-10-04"
Set-TextCell $ws.Range("G27") "2023-10-13"
$ws.Range("H27").Value = 8000
$ws.Range("I27").Value = 4000000
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 100

